$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the "current resource" row that feeds the formulas in B4:B6 ---
$ws.Range("A2").Value = "211_NE_MHLC"
$ws.Range("B2").Value = "211 Northeast Region"
$ws.Range("C2").Value = "MHLCdata()"
$ws.Range("D2").Value = "211 Northeast Region"
$ws.Range("E2").Value = "T"

# Match the new formatting (font/number-format) applied to A2, C2 and D2 by
# copying the format already used elsewhere in the sheet (style index 5,
# e.g. I9) onto each of them individually (not the whole row, so B2 keeps
# its original formatting).
$fmtSrc = $ws.Range("I9")

$fmtSrc.Copy()
$ws.Range("A2").PasteSpecial(-4122)

$fmtSrc.Copy()
$ws.Range("C2").PasteSpecial(-4122)

$fmtSrc.Copy()
$ws.Range("D2").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Rows 10-12: replace the leftover empty-formula rows with the new
#     "211 Northeast Region" module code (static text, like the author
#     pasted-as-values after generating it once). ---
$ws.Range("B10").Value = "mod_Accordion_ui('211_NE')"
$ws.Range("B11").Value = "mod_Accordion_server('211_NE', selector=selection, data=FHFdata(), title = c('211 Northeast Region'), Visible = T)"
$ws.Range("B12").Value = "mod_info_server('211_NE', selector = selection, data = FHFdata(), rownametitle = c('211 Northeast Region'), phone = T, website = T, email = F)"

# --- Move/record the active selection like the author left it ---
$ws.Range("B4").Select()
